$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $c = $ws.Range($cellAddr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '61.445.70'
Set-TextValue "E2" '  +0.57%  '

Set-TextValue "D3" '2.934.37'

Set-TextValue "D4" '1.00'
Set-TextValue "E4" '  +0.00%  '

Set-TextValue "D5" '595.44'
Set-TextValue "E5" '  +0.85%  '

Set-TextValue "D6" '145.16'
Set-TextValue "E6" '  -0.05%  '

Set-TextValue "E7" '  +0.02%  '

Set-TextValue "D8" '0.503'
Set-TextValue "E8" '  -0.66%  '

Set-TextValue "D9" '6.96'
Set-TextValue "E9" '  +1.84%  '

Set-TextValue "E10" '  -1.52%  '

Set-TextValue "D11" '0.440'
Set-TextValue "E11" '  -0.56%  '

Set-TextValue "E12" '  -0.60%  '

Set-TextValue "D13" '33.69'
Set-TextValue "E13" '  -0.03%  '

Set-TextValue "E14" '  +0.62%  '

Set-TextValue "D15" '3.421.08'
Set-TextValue "E15" '  +0.43%  '

Set-TextValue "D16" '61.427.76'
Set-TextValue "E16" '  +0.59%  '

Set-TextValue "E17" '  +0.24%  '

Set-TextValue "D18" '2.932.27'

Set-TextValue "D19" '432.48'
Set-TextValue "E19" '  +0.31%  '

Set-TextValue "D20" '13.50'
Set-TextValue "E20" '  +0.19%  '

Set-TextValue "E21" '  -0.57%  '

Set-TextValue "D22" '7.13'
Set-TextValue "E22" '  +0.57%  '

Set-TextValue "D23" '81.89'
Set-TextValue "E23" '  +1.05%  '

Set-TextValue "D24" '10.87'
Set-TextValue "E24" '  -1.12%  '

Set-TextValue "E25" '  -1.14%  '

Set-TextValue "D26" '11.81'
Set-TextValue "E26" '  -1.96%  '

Set-TextValue "E27" '  -0.04%  '

Set-TextValue "D28" '2.23'
Set-TextValue "E28" '  -3.21%  '

Set-TextValue "E29" '  -0.21%  '

Set-TextValue "E30" '  -2.72%  '

Set-TextValue "E31" '  +1.64%  '

Set-TextValue "D32" '26.72'
Set-TextValue "E32" '  +0.65%  '

Set-TextValue "E33" '  +0.04%  '

Set-TextValue "E34" '  +2.16%  '

Set-TextValue "E35" '  -0.21%  '

Set-TextValue "E36" '  +0.17%  '

Set-TextValue "E37" '  -1.84%  '

Set-TextValue "E38" '  +0.37%  '

Set-TextValue "E39" '  -0.66%  '

Set-TextValue "D40" '8.62'
Set-TextValue "E40" '  +0.05%  '

Set-TextValue "D41" '42.77'
Set-TextValue "E41" '  +8.82%  '

Set-TextValue "E42" '  -1.40%  '

Set-TextValue "E43" '  -0.15%  '

Set-TextValue "D44" '2.706.97'
Set-TextValue "E44" '  -0.16%  '

Set-TextValue "B45" 'Monero'
Set-TextValue "C45" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D45" '134.27'
Set-TextValue "E45" '  +2.01%  '

Set-TextValue "B46" 'Bittensor'
Set-TextValue "C46" 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue "D46" '365.21'
Set-TextValue "E46" '  -3.13%  '

Set-TextValue "E47" '  +0.01%  '

Set-TextValue "E48" '  -1.82%  '

Set-TextValue "E49" '  -1.30%  '

Set-TextValue "D50" '2.00'
Set-TextValue "E50" '  -1.42%  '

Set-TextValue "E51" '  -0.48%  '
